$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Season 16 was mistakenly split across two rows (row 17 "16 a" / row 18 "16 b").
# Merge row 18's episode-number cells into row 17 (continuing the episode
# count sequence), fix the season number, then remove the now-redundant row 18.

# 1) Copy row 18's episode cells (C18:W18 - 21 columns, values 1..20 plus the
#    season-episode-count marker in column N) into row 17 right after its
#    existing episode cells (which end at AA17, i.e. 25 episodes).
$ws.Range("C18:W18").Copy()
$ws.Range("AB17").PasteSpecial()

# 2) The pasted episode numbers need to continue the season's running count
#    (row 17 already holds episodes 1-25), so add 25 to every pasted episode
#    number EXCEPT the season-episode-count marker cell (originally N18,
#    now AM17) which must keep its original value of 16.
$offsetCols = @("AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AN","AO","AP","AQ","AR","AS","AT","AU","AV")
foreach ($col in $offsetCols) {
    $cell = $ws.Range($col + "17")
    $cell.Value = $cell.Value() + 25
}

# 3) Row 17's "S." column held the text label "16 a"; season 16 is a single
#    season now, so it becomes the plain number 16 (matching every other row).
$ws.Range("A17").Value = 16

# 4) Row 18 (the old "16 b" / "Adventures in Unova and Beyond" row) is now
#    fully folded into row 17, so delete it - this shifts every row below
#    up by one, which is exactly what the renumbered season table needs.
$ws.Rows(18).Delete()
